# Update "想去人数" (interested count) values for several events that
# appear on both the "展览" sheet and the consolidated "全部类型" sheet.
# Source: gh-pages data regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - F column holds counts, rows keyed by row number
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 303
$wsExpo.Range("F4").Value = 3689
$wsExpo.Range("F5").Value = 2253
$wsExpo.Range("F6").Value = 439
$wsExpo.Range("F8").Value = 8
$wsExpo.Range("F12").Value = 1366
$wsExpo.Range("F14").Value = 2103
$wsExpo.Range("F15").Value = 153

# Sheet "全部类型" (All Types) - same events, different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 303
$wsAll.Range("F4").Value = 3689
$wsAll.Range("F5").Value = 2253
$wsAll.Range("F6").Value = 439
$wsAll.Range("F8").Value = 8
$wsAll.Range("F15").Value = 1366
$wsAll.Range("F17").Value = 2103
$wsAll.Range("F18").Value = 153
